$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("B2").Value = 3.286832544864788
$ws.Range("C2").Value = 3286.919754855326
$ws.Range("D2").Value = 6708.013860684405
$ws.Range("E2").Value = 1133.036916526867
$ws.Range("G2").Value = 11131.25736461146

# Row 3
$ws.Range("B3").Value = 3.286832544864788
$ws.Range("C3").Value = 1.655778082260271
$ws.Range("D3").Value = 0.7527432677738641
$ws.Range("E3").Value = 0.4942365360607697
$ws.Range("G3").Value = 6.189590430959694

# Row 4
$ws.Range("B4").Value = 0.6606524410359556
$ws.Range("C4").Value = 1.655778082260271
$ws.Range("D4").Value = 261.3203778131603
$ws.Range("E4").Value = 10.19245300693656
$ws.Range("G4").Value = 273.8292613433931
